# Populate Sheet1 with the semivariogram model header row + first data row
# (state/year/model/c0/c0_c1/a/gde/rss/r2 for PA, 2015, Exp model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1 already carries the bold/centered style in the template)
$headers = @("state", "year", "model", "c0", "c0_c1", "a", "gde", "rss", "r2")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# First data row
$ws.Cells.Item(2, 1).Value = "PA"
$ws.Cells.Item(2, 2).Value = 2015
$ws.Cells.Item(2, 3).Value = "Exp"
$ws.Cells.Item(2, 4).Value = 242.7104
$ws.Cells.Item(2, 5).Value = 665.8231
$ws.Cells.Item(2, 6).Value = 282.15
$ws.Cells.Item(2, 7).Value = 0.3645268540547782
$ws.Cells.Item(2, 8).Value = 159021056.9176
$ws.Cells.Item(2, 9).Value = -5.6901
